$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.890575333333333
$ws.Range("H2").Value = 5.671726
$ws.Range("I2").Value = 0.006346320422088561
$ws.Range("J2").Value = 0.00634632042208856
$ws.Range("M2").Value = 2.750415333333333
$ws.Range("N2").Value = 8.251245999999998
$ws.Range("O2").Value = 0.04811444325525444
$ws.Range("P2").Value = 0.04811444325525444
$ws.Range("Q2").Value = 5.199867385621777
$ws.Range("R2").Value = 46.79880647059598
$ws.Range("S2").Value = 0.0003053496738282425
$ws.Range("T2").Value = 0.0003053496738282425
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.890575333333333
$ws.Range("H3").Value = 5.671726
$ws.Range("I3").Value = 0.006346320422088561
$ws.Range("J3").Value = 0.00634632042208856
$ws.Range("O3").Value = 0.01415294505639593
$ws.Range("P3").Value = 0.01415294505639593
$ws.Range("Q3").Value = 1.529549807296444
$ws.Range("R3").Value = 13.765948265668
$ws.Range("S3").Value = 0.00008981912424410281
$ws.Range("T3").Value = 0.00008981912424410281
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.890575333333333
$ws.Range("H4").Value = 5.671726
$ws.Range("I4").Value = 0.006346320422088561
$ws.Range("J4").Value = 0.00634632042208856
$ws.Range("O4").Value = 0.9377326116883496
$ws.Range("P4").Value = 0.9377326116883496
$ws.Range("Q4").Value = 101.3434822072824
$ws.Range("R4").Value = 912.091339865542
$ws.Range("S4").Value = 0.005951151624016215
$ws.Range("T4").Value = 0.005951151624016215
$ws.Range("I5").Value = 0.8887896079640043
$ws.Range("J5").Value = 0.8887896079640044
$ws.Range("M5").Value = 2.750415333333333
$ws.Range("N5").Value = 8.251245999999998
$ws.Range("O5").Value = 0.04811444325525444
$ws.Range("P5").Value = 0.04811444325525444
$ws.Range("Q5").Value = 728.2311304430854
$ws.Range("R5").Value = 6554.080173987769
$ws.Range("S5").Value = 0.04276361715824393
$ws.Range("T5").Value = 0.04276361715824394
$ws.Range("I6").Value = 0.8887896079640043
$ws.Range("J6").Value = 0.8887896079640044
$ws.Range("O6").Value = 0.01415294505639593
$ws.Range("P6").Value = 0.01415294505639593
$ws.Range("S6").Value = 0.01257899048821023
$ws.Range("T6").Value = 0.01257899048821023
$ws.Range("I7").Value = 0.8887896079640043
$ws.Range("J7").Value = 0.8887896079640044
$ws.Range("O7").Value = 0.9377326116883496
$ws.Range("P7").Value = 0.9377326116883496
$ws.Range("S7").Value = 0.8334470003175501
$ws.Range("T7").Value = 0.8334470003175503
$ws.Range("I8").Value = 0.104864071613907
$ws.Range("J8").Value = 0.104864071613907
$ws.Range("M8").Value = 2.750415333333333
$ws.Range("N8").Value = 8.251245999999998
$ws.Range("O8").Value = 0.04811444325525444
$ws.Range("P8").Value = 0.04811444325525444
$ws.Range("Q8").Value = 85.92053814534799
$ws.Range("R8").Value = 773.2848433081318
$ws.Range("S8").Value = 0.005045476423182265
$ws.Range("T8").Value = 0.005045476423182265
$ws.Range("I9").Value = 0.104864071613907
$ws.Range("J9").Value = 0.104864071613907
$ws.Range("O9").Value = 0.01415294505639593
$ws.Range("P9").Value = 0.01415294505639593
$ws.Range("S9").Value = 0.001484135443941593
$ws.Range("T9").Value = 0.001484135443941593
$ws.Range("I10").Value = 0.104864071613907
$ws.Range("J10").Value = 0.104864071613907
$ws.Range("O10").Value = 0.9377326116883496
$ws.Range("P10").Value = 0.9377326116883496
$ws.Range("S10").Value = 0.0983344597467831
$ws.Range("T10").Value = 0.0983344597467831
